# Practica 3 Unidad 2
# Apply changes to the "Inventario" worksheet (sheet2):
#  - remove the "Bio-Electro 24 tabletas" and "Almetec 40mg 28 tabletas" rows
#  - update several "Stock" quantities
#  - change "Talco para bebé" entry to "Talco para pies" with new price/stock
#  - append a new product "Aspirina 250mg" at the end

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventario")

# Remove the two rows (Bio-Electro, Almetec) - row 12 first so row numbers don't shift
$ws.Rows.Item(12).Delete() | Out-Null
$ws.Rows.Item(11).Delete() | Out-Null

# Update Stock (column D) values
$ws.Range("D2").Value = 21
$ws.Range("D3").Value = 0
$ws.Range("D5").Value = 23
$ws.Range("D9").Value = 11
$ws.Range("D11").Value = 17
$ws.Range("D13").Value = 43
$ws.Range("D14").Value = 5

# Replace the "Talco para bebé" row with "Talco para pies" details
$ws.Range("B7").Value = "Talco para pies"
$ws.Range("C7").Value = 25
$ws.Range("D7").Value = 5

# Append the new product row (Folio stored as text, like the other Folio values)
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "6771"
$ws.Range("B16").Value = "Aspirina 250mg"
$ws.Range("C16").Value = 35
$ws.Range("D16").Value = 20
